# Fruta / hortaliza, semanal
# The weekly refresh reshuffles the per-record fields (date, quality,
# volume, price bounds, unit of sale, origin price and kg/unit) across
# the existing data rows (2-16). Columns A,B,C,E,F,G,H,I,J,K,R are
# constant across all rows already, so only D,L,M,N,O,P,Q,S,T move.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("D","L","M","N","O","P","Q","S","T")

# Snapshot every source row's values BEFORE writing anything, since the
# row->row remapping below contains cycles (e.g. 2<->9) and longer
# chains, so in-place writes would clobber data still needed later.
$snapshot = @{}
for ($r = 2; $r -le 16; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value()
    }
    $snapshot[$r] = $rowVals
}

# target row -> source row (the row whose D/L/M/N/O/P/Q/S/T values now
# land on the target row)
$mapping = @{
    2  = 9
    3  = 15
    4  = 16
    5  = 12
    6  = 11
    7  = 6
    8  = 5
    9  = 2
    10 = 8
    11 = 4
    12 = 13
    13 = 7
    14 = 14
    15 = 10
    16 = 3
}

for ($target = 2; $target -le 16; $target++) {
    $source = $mapping[$target]
    $srcVals = $snapshot[$source]
    foreach ($c in $cols) {
        $ws.Range("$c$target").Value = $srcVals[$c]
    }
}
